$d = $word.ActiveDocument

# The last paragraph currently is "Hasta la vista baby! (Terminator 2)"
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# Insert a brand new paragraph right after it (inherits formatting, incl. pt-BR lang)
$lastRange.InsertParagraphAfter()

# Move into the newly created (now last) paragraph and type the new quote
$newRange = $d.Paragraphs.Last.Range
$newRange.Collapse(1)  # wdCollapseStart
$newRange.InsertAfter("Hello World!")
